$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range($ws2.Cells.Item(1,5), $ws2.Cells.Item(1,16384)).EntireColumn.ColumnWidth = 11.43
